$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 10.01.2022 08:04"

# Row 10: new price scraped, old price shifted into "Old Cena", and the
# delta / old-datum columns got written as plain text (as produced by the
# AWS bash cmd line scraper) instead of as numeric / date values.
$ws.Range("B10").Value = 36.2
$ws.Range("C10").Value = 36.5

$ws.Range("D10").Value = "'-0.3"
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = "'2022-01-10 08:04:35"
$ws.Range("E10").Style = "Normal"
